# 9.3.1 data sheet update:
#  - add a new "2021" column (O) mirroring the formatting of the existing
#    "2020" column (N) for the header block (row 3 bottom border), the
#    year-label row (row 4) and the data row (row 5)
#  - revise a couple of previously-entered data points in row 5
#  - move the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend column O so it inherits the same look as column N ----------

# row 3 is the thin bottom-border spacer row above the year headers
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)   # xlPasteFormats

# row 4 holds the year labels (2010 ... 2020) -> add 2021
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O4").Value = 2021

# row 5 holds the data values for each year -> add the 2021 figure
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O5").Value = 4.0999999999999996

# --- revise existing data points in row 5 -------------------------------
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1

# --- move the selected cell ---------------------------------------------
$ws.Range("P4").Select() | Out-Null
